$d = $word.ActiveDocument
$bullet = [char]0x2022

# --- 1. Collapse CORE COMPETENCIES from three long bullet paragraphs down to one short line ---
# Paragraph 6 = "Survey Methodology & Research Design: ..."
# Paragraph 7 = "Redistricting & Geospatial Analysis: ..."
# Paragraph 8 = "Data Analysis & Visualization: ..."
# Delete paragraphs 8 and 7 (highest index first so indices of earlier paragraphs stay valid),
# then overwrite paragraph 6's text with the condensed summary line.
$d.Paragraphs(8).Range.Delete()
$d.Paragraphs(7).Range.Delete()

$p6 = $d.Paragraphs(6)
$p6.Range.Text = "Survey Methodology & Research Design " + $bullet + " Redistricting & Geospatial Analysis " + $bullet + " Data Analysis & Visualization"

# --- 2. Append a new "TECHNICAL SKILLS" section at the end of the document with the detailed content ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)

# Heading paragraph
$lastPara.Range.InsertParagraphAfter()
$count = $d.Paragraphs.Count
$headingPara = $d.Paragraphs($count)
$headingPara.Style = "Heading 2"
$headingPara.Range.Text = "TECHNICAL SKILLS"

# Body paragraph 1
$headingPara.Range.InsertParagraphAfter()
$count = $d.Paragraphs.Count
$body1 = $d.Paragraphs($count)
$body1.Style = "Normal"
$body1.Range.Text = "SURVEY METHODOLOGY & RESEARCH DESIGN Survey Design and Questionnaire Development for Political and Market Research; Sampling Methodology and Statistical Analysis (R, SPSS, Stata, OSCAR); Random Device Engagement (RDE), Text Message, Web Panel, and Live Telephone Calling; Focus Groups and Qualitative Research Methodologies; Meta-analytical Dataset Development for Longitudinal Analysis; Survey Instrument Standardization and Call Methods Optimization; Expert Testimony and Consultation on Research Methodology"

# Body paragraph 2
$body1.Range.InsertParagraphAfter()
$count = $d.Paragraphs.Count
$body2 = $d.Paragraphs($count)
$body2.Style = "Normal"
$body2.Range.Text = "REDISTRICTING & GEOSPATIAL ANALYSIS Redistricting Software Development and Boundary Estimation Systems; Geospatial Analysis; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Spatial Clustering and Boundary Estimation without ML Requirements; Census Data Integration and Demographic Mapping; Court Case Analysis and Expert Testimony for Redistricting; Multi-tenant Data Warehouse Design for Electoral Analytics"

# Body paragraph 3
$body2.Range.InsertParagraphAfter()
$count = $d.Paragraphs.Count
$body3 = $d.Paragraphs($count)
$body3.Style = "Normal"
$body3.Range.Text = "DATA ANALYSIS & VISUALIZATION Advanced Statistical Modeling and Analysis (Regression, Clustering, Segmentation); Data Visualization; Consumer Behavior Analysis and Market Segmentation; Machine Learning and Predictive Modeling for Targeting; Big Data Analytics; Fraud Detection and Entity Resolution Systems; Multi-million Dollar Research Project Management"

Write-Output "Done. Final paragraph count: $($d.Paragraphs.Count)"
